# Generate Report for Handoff
# The file "4821b8ad-595a-4a63-95d6-097abe830311" has moved from
# "Handed back: in sync with en-US" to "Ready for handoff" for both
# the zh-cn and de-de target languages. Update the per-language sheets
# and roll the status up into the Overview sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("E3").Value = "2016-03-21 16:49:29"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("E3").Value = "2016-03-21 16:49:33"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = "2016-03-21 16:49:33"
